# Rename transcript speaker "RBD" to "T" in column D of the DataSheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 4)  # Column D
    if ($cell.Value2 -eq "RBD") {
        $cell.Value = "T"
    }
}
